# Weekly crime-data refresh: bump the report volume/number + date range in
# the header, and refresh the Week/28-day/YTD/2yr crime-count figures (and
# their derived % change columns) for several precincts' categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 30   Number  13" -> "...Number  14"
# and     "Report Covering the Week  3/27/2023  Through  4/2/2023"
#      -> "...4/3/2023  Through  4/9/2023"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# ---------------------------------------------------------------------
# Cells that flip from a text placeholder ("0" / "***.*") to a real
# number, or vice versa, need their number format copied from a sibling
# cell in the same row so the stored style index matches (General-text
# style "14" for placeholders, "#,##0" style "16" for counts, "#,##0.0"
# style "15" for percentages).
# ---------------------------------------------------------------------

# C15: "0" (text) -> 1 (number)
$ws.Range("C15").Value = 1
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# G15: 1 (number) -> "0" (text)
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# H15: 100 (number) -> "***.*" (text)
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C18: "0" (text) -> 1 (number)
$ws.Range("C18").Value = 1
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C26: "0" (text) -> 2 (number) ; D26: "0" (text) -> 1 (number)
$ws.Range("C26").Value = 2
$ws.Range("F26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D26").Value = 1
$ws.Range("F26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# E26: "***.*" (text) -> 100 (number)
$ws.Range("E26").Value = 100
$ws.Range("H26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C27: "0" (text) -> 1 (number)
$ws.Range("C27").Value = 1
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Plain numeric refreshes (style unchanged)
# ---------------------------------------------------------------------

# Row 15
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -25

# Row 16
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 59
$ws.Range("J16").Value = 56
$ws.Range("K16").Value = 5.357142857142
$ws.Range("L16").Value = 63.888888888888
$ws.Range("M16").Value = -19.178082191780
$ws.Range("N16").Value = -78.776978417266

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 4.761904761904
$ws.Range("I17").Value = 77
$ws.Range("J17").Value = 72
$ws.Range("K17").Value = 6.944444444444
$ws.Range("L17").Value = 16.666666666666
$ws.Range("M17").Value = 133.333333333333
$ws.Range("N17").Value = -14.444444444444

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -80
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = -31.578947368421
$ws.Range("L18").Value = 4
$ws.Range("M18").Value = -64.864864864864
$ws.Range("N18").Value = -92.419825072886

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("E19").Value = -21.428571428571
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -13.725490196078
$ws.Range("I19").Value = 149
$ws.Range("J19").Value = 206
$ws.Range("K19").Value = -27.669902912621
$ws.Range("L19").Value = 28.448275862069
$ws.Range("M19").Value = 93.506493506493
$ws.Range("N19").Value = -1.973684210526

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -5.882352941176
$ws.Range("I20").Value = 70
$ws.Range("J20").Value = 84
$ws.Range("K20").Value = -16.666666666666
$ws.Range("L20").Value = 66.666666666666
$ws.Range("M20").Value = -5.405405405405
$ws.Range("N20").Value = -92.592592592592

# Row 21
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 115
$ws.Range("H21").Value = -10.434782608695
$ws.Range("I21").Value = 387
$ws.Range("J21").Value = 462
$ws.Range("K21").Value = -16.233766233766
$ws.Range("L21").Value = 33.448275862069
$ws.Range("M21").Value = 13.823529411764
$ws.Range("N21").Value = -78.817733990147

# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -41.935483870967
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -7.575757575757
$ws.Range("I24").Value = 362
$ws.Range("J24").Value = 365
$ws.Range("K24").Value = -0.821917808219
$ws.Range("L24").Value = 69.158878504672
$ws.Range("M24").Value = 118.072289156627

# Row 25
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 27.272727272727
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 32.352941176470
$ws.Range("I25").Value = 140
$ws.Range("J25").Value = 106
$ws.Range("K25").Value = 32.075471698113
$ws.Range("L25").Value = 34.615384615384
$ws.Range("M25").Value = 3.703703703703

# Row 26
$ws.Range("F26").Value = 5
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -16.666666666666
$ws.Range("L26").Value = 42.857142857142

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("I27").Value = 19
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = 5.555555555555
$ws.Range("L27").Value = 72.727272727272

# Row 28 / 29 (2-year % change only)
$ws.Range("N28").Value = -76.470588235294
$ws.Range("N29").Value = -88.235294117647
